$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("sistemas")
$dst = $wb.Worksheets.Item("usuarios")
$src.Rows.Item(2).Copy()
$dst.Rows.Item(2).PasteSpecial(-4122)
